# Apply the edit: insert one new data row before the current row 64
# (shifting the existing rows 64..116 down to 65..117), then populate
# the newly inserted row 64 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 64; this shifts rows 64..116 down to 65..117
# and the cells in the new row 64 are blank.
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the new record.
$ws.Cells.Item(64, 1).Value = 5
$ws.Cells.Item(64, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(64, 3).Value = 'Maule'
$ws.Cells.Item(64, 4).Value = 45216
$ws.Cells.Item(64, 5).Value = 7
$ws.Cells.Item(64, 6).Value = 300000000
$ws.Cells.Item(64, 7).Value = 'Espárragos'
$ws.Cells.Item(64, 8).Value = 'Verde'
$ws.Cells.Item(64, 9).Value = 'Primera'
$ws.Cells.Item(64, 10).Value = 4000
$ws.Cells.Item(64, 11).Value = 1100
$ws.Cells.Item(64, 12).Value = 1200
$ws.Cells.Item(64, 13).Value = 1150
$ws.Cells.Item(64, 14).Value = '$/kilo'
$ws.Cells.Item(64, 15).Value = 'Provincia de Linares'
$ws.Cells.Item(64, 16).Value = 1150
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = 'Hortaliza'

# Ensure the date cell keeps the same number format as the rest of
# column D (it should inherit it from the row-insert, but set it
# explicitly to be safe).
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
